# Updated cryptos list on Tue Apr  2 23:09:45 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns of the crypto table,
# and re-syncs two pairs of rows whose ranking order changed (Stacks/PEPE
# at rows 39-40, InjectiveProtocol/CoreDAO at rows 43-44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Force the literal string into the cell as TEXT (matches the source
    # workbook, where every Price cell is stored as an inline string) even
    # when the new value happens to look like a plain number (e.g. "1.00").
    # A leading apostrophe is Excel's "treat as text" quote-prefix; we then
    # reset the cell style back to Normal so no stray number-format/
    # quote-prefix style is left behind on the cell.
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

# --- Price (D) / Volume 1h (E) updates ---------------------------------
Set-TextValue "D2"  "65.792.89"
$ws.Range("E2").Value  = "  -5.76%  "

Set-TextValue "D3"  "3.295.25"
$ws.Range("E3").Value  = "  -6.06%  "

$ws.Range("E4").Value  = "  +0.06%  "

Set-TextValue "D5"  "557.47"
$ws.Range("E5").Value  = "  -3.54%  "

Set-TextValue "D6"  "183.66"
$ws.Range("E6").Value  = "  -4.89%  "

Set-TextValue "D7"  "1.00"
$ws.Range("E7").Value  = "  +0.06%  "

$ws.Range("E8").Value  = "  -3.94%  "

Set-TextValue "D9"  "3.291.90"
$ws.Range("E9").Value  = "  -5.86%  "

$ws.Range("E10").Value = "  -9.93%  "

Set-TextValue "D11" "0.584"
$ws.Range("E11").Value = "  -6.01%  "

Set-TextValue "D12" "47.34"
$ws.Range("E12").Value = "  -8.13%  "

Set-TextValue "D13" "0.0000266"
$ws.Range("E13").Value = "  -7.06%  "

Set-TextValue "D14" "644.07"
$ws.Range("E14").Value = "  -0.32%  "

Set-TextValue "D15" "8.64"
$ws.Range("E15").Value = "  -5.85%  "

Set-TextValue "D16" "3.823.16"
$ws.Range("E16").Value = "  -5.98%  "

Set-TextValue "D17" "18.10"
$ws.Range("E17").Value = "  -1.60%  "

Set-TextValue "D18" "65.784.12"
$ws.Range("E18").Value = "  -5.76%  "

$ws.Range("E19").Value = "  -3.19%  "

Set-TextValue "D20" "3.294.01"
$ws.Range("E20").Value = "  -6.11%  "

Set-TextValue "D21" "11.39"
$ws.Range("E21").Value = "  -8.27%  "

$ws.Range("E22").Value = "  -4.86%  "

Set-TextValue "D23" "18.31"
$ws.Range("E23").Value = "  +0.69%  "

Set-TextValue "D24" "107.74"
$ws.Range("E24").Value = "  +8.57%  "

$ws.Range("E25").Value = "  -7.92%  "

$ws.Range("E26").Value = "  -7.37%  "

$ws.Range("E27").Value = "  -7.10%  "

Set-TextValue "D28" "9.58"
$ws.Range("E28").Value = "  -5.38%  "

Set-TextValue "D29" "8.68"
$ws.Range("E29").Value = "  -7.32%  "

Set-TextValue "D30" "30.31"
$ws.Range("E30").Value = "  -7.43%  "

Set-TextValue "D31" "3.91"
$ws.Range("E31").Value = "  -8.62%  "

$ws.Range("E32").Value = "  -6.76%  "

Set-TextValue "D33" "11.07"
$ws.Range("E33").Value = "  -5.09%  "

$ws.Range("E34").Value = "  -4.67%  "

Set-TextValue "D35" "3.803.23"
$ws.Range("E35").Value = "  +1.43%  "

Set-TextValue "D36" "57.49"
$ws.Range("E36").Value = "  -6.66%  "

$ws.Range("E37").Value = "  -0.02%  "

Set-TextValue "D38" "520.65"
$ws.Range("E38").Value = "  -8.38%  "

# --- Rows 39/40: Stacks and PEPE swap ranking positions -----------------
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D39" "0.0₃0737"
$ws.Range("E39").Value = "  -7.23%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D40" "3.38"
$ws.Range("E40").Value = "  -6.96%  "

$ws.Range("E41").Value = "  -2.22%  "

$ws.Range("E42").Value = "  -6.18%  "

# --- Rows 43/44: InjectiveProtocol and CoreDAO swap ranking positions ---
$ws.Range("B43").Value = "CoreDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue "D43" "3.38"
$ws.Range("E43").Value = "  -13.28%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D44" "32.98"
$ws.Range("E44").Value = "  -4.04%  "

Set-TextValue "D45" "0.338"
$ws.Range("E45").Value = "  -9.91%  "

Set-TextValue "D46" "0.0414"
$ws.Range("E46").Value = "  -6.61%  "

Set-TextValue "D47" "3.25"
$ws.Range("E47").Value = "  -2.52%  "

Set-TextValue "D48" "0.129"
$ws.Range("E48").Value = "  -4.59%  "

$ws.Range("E49").Value = "  -9.39%  "

Set-TextValue "D50" "0.999"
$ws.Range("E50").Value = "  +0.16%  "

Set-TextValue "D51" "1.26"
$ws.Range("E51").Value = "  +1.13%  "
